# Update automatico via Actualizar 02-22-2021 13-01-07
#
# This sheet keeps a rolling log of "disponibilidad" refresh timestamps in
# column D, grouped in blocks of 14 rows (one block per room/category set).
# Each time the source is refreshed, the newest timestamp block is pushed
# down to the next older block, and a brand-new timestamp is recorded for
# the most-recent block. Concretely, for this run:
#   - rows 2:15  (newest block) get the new refresh timestamp
#   - rows 16:29 get what used to be the previous "newest" timestamp
#   - rows 30:43 get what used to be the middle timestamp
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:D15").Value = 44249.54226155607
$ws.Range("D16:D29").Value = 44249.52093606482
$ws.Range("D30:D43").Value = 44249.49960440972
